$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the performance table (A1:D34) by column D (Accuracy) descending,
# instead of the previous sort by column C (AUC).
$dataRange = $ws.Range("A1:D34")
$sortKey = $ws.Range("D1:D34")
[void]$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)

# Move the selection to D6, matching where the user ended up after sorting.
[void]$ws.Range("D6").Select()
